# Updates cryptos list values (Price / Volume(1h) columns, plus a few
# coin-name/link swaps) to match the latest scraped snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.280.04'
$ws.Range("E2").Value = '  -1.05%  '
$ws.Range("D3").Value = '1.848.73'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.33%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.52'
$ws.Range("E5").Value = '  -2.94%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.005'
$ws.Range("E6").Value = '  +0.34%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4646'
$ws.Range("E7").Value = '  -0.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3868'
$ws.Range("E8").Value = '  -0.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07870'
$ws.Range("E9").Value = '  -0.27%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9650'
$ws.Range("E10").Value = '  -1.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.10'
$ws.Range("E11").Value = '  -0.61%  '
$ws.Range("D12").Value = '1.838.16'
$ws.Range("E12").Value = '  +0.70%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.696'
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.898'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06933'
$ws.Range("E15").Value = '  +0.48%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.80'
$ws.Range("E16").Value = '  +1.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.006'
$ws.Range("E17").Value = '  +0.30%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009987'
$ws.Range("E18").Value = '  -0.15%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.76'
$ws.Range("E19").Value = '  -1.66%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.006'
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("D21").Value = '28.264.00'
$ws.Range("E21").Value = '  -1.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.316'
$ws.Range("E22").Value = '  -1.17%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.05'
$ws.Range("E23").Value = '  -1.46%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.109'
$ws.Range("E24").Value = '  -3.19%  '
$ws.Range("D25").Value = '2.046.09'
$ws.Range("E25").Value = '  -0.55%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.68'
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.22'
$ws.Range("E27").Value = '  -0.91%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.749'
$ws.Range("E28").Value = '  -4.22%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.976'
$ws.Range("E29").Value = '  -1.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.04'
$ws.Range("E30").Value = '  +1.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09280'
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9311'
$ws.Range("E32").Value = '  -3.83%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.290'
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.327'
$ws.Range("E34").Value = '  -1.16%  '
$ws.Range("E35").Value = '  -3.71%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05828'
$ws.Range("E36").Value = '  -4.19%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02109'
$ws.Range("E37").Value = '  -3.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.143'
$ws.Range("E38").Value = '  -1.72%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.799'
$ws.Range("E39").Value = '  +1.96%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5617'
$ws.Range("E40").Value = '  -1.22%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.937'
$ws.Range("E41").Value = '  -1.98%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1767'
$ws.Range("E42").Value = '  -1.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07202'
$ws.Range("E43").Value = '  +1.57%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.67'
$ws.Range("E44").Value = '  -0.33%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5286'
$ws.Range("E45").Value = '  -1.48%  '
$ws.Range("B46").Value = 'RenderToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.144'
$ws.Range("E46").Value = '  -10.86%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.136'
$ws.Range("E47").Value = '  -7.23%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.840'
$ws.Range("E48").Value = '  -3.14%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.13'
$ws.Range("E49").Value = '  +0.42%  '
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("B51").Value = 'EOS'
$ws.Range("C51").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.025'
$ws.Range("E51").Value = '  +0.39%  '
